# "fixed issue of yasumi day"
# The lesson deck skipped a holiday class, so the last 5 "placeholder/future"
# slides (10-14) are removed, and the remaining 9 kanji slides are updated to
# the next set of vocabulary (with the page reference moving from 69-70 to
# 67-68).

$p = $ppt.ActivePresentation

# --- Remove the trailing 5 slides (old slides 10-14) -----------------------
for ($i = $p.Slides.Count; $i -ge 10; $i--) {
    $p.Slides.Item($i).Delete()
}

# --- Update remaining 9 slides with the new kanji content -------------------
# Each slide has 4 shapes: 1=kanji, 2=reading, 3=definition, 4=page range.
# The reading shape has 2 leading blank paragraphs before the actual text,
# so it must be edited via Paragraphs(3) to keep that layout intact.

function Set-KanjiSlide {
    param($Slide, $Kanji, $Reading, $Definition, $Pages)

    $Slide.Shapes.Item(1).TextFrame.TextRange.Text = $Kanji
    $Slide.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Text = $Reading
    $Slide.Shapes.Item(3).TextFrame.TextRange.Text = $Definition
    $Slide.Shapes.Item(4).TextFrame.TextRange.Text = $Pages
}

Set-KanjiSlide $p.Slides.Item(1) "殺す" "ころす" "to kill, to slay, to murder, to slaughter | to suppress, to block, to hamper, to destroy (e.g. talent), to eliminate (e.g..." "67-68"

Set-KanjiSlide $p.Slides.Item(2) "殺人" "さつじん" "murder, homicide, manslaughter..." "67-68"

Set-KanjiSlide $p.Slides.Item(3) "農薬" "のうやく" "agricultural chemical (i.e. pesticide, herbicide, fungicide, etc.), agrochemical, agrichemical..." "67-68"

Set-KanjiSlide $p.Slides.Item(4) "収入印紙" "しゅうにゅういんし" "revenue stamp..." "67-68"

Set-KanjiSlide $p.Slides.Item(5) "収穫" "しゅうかく" "harvest, crop, ingathering | fruits (of one's labors), gain, result, returns..." "67-68"

Set-KanjiSlide $p.Slides.Item(6) "少量" "しょうりょう" "small quantity, small amount | narrowmindedness..." "67-68"

Set-KanjiSlide $p.Slides.Item(7) "完全" "かんぜん" "perfect, complete..." "67-68"

Set-KanjiSlide $p.Slides.Item(8) "原因" "げんいん" "cause, origin, source..." "67-68"

Set-KanjiSlide $p.Slides.Item(9) "一環" "いっかん" "link (e.g. in a chain of events), part (of a plan, campaign, activities, etc.) | monocyclic..." "67-68"
